$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs; all runs share identical formatting) ---
$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("M14").Value = -57.142857142857
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 35.714285714285
$ws.Range("L15").Value = 171.428571428571
$ws.Range("M15").Value = 72.727272727272
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -39.285714285714
$ws.Range("I16").Value = 133
$ws.Range("J16").Value = 142
$ws.Range("K16").Value = -6.338028169014
$ws.Range("L16").Value = 12.711864406779
$ws.Range("M16").Value = -10.738255033557
$ws.Range("N16").Value = -64.343163538874
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 67.647058823529
$ws.Range("I17").Value = 257
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 16.818181818181
$ws.Range("L17").Value = 43.575418994413
$ws.Range("M17").Value = 90.37037037037
$ws.Range("N17").Value = 54.819277108433
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -56.521739130434
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 127
$ws.Range("K18").Value = -25.984251968503
$ws.Range("L18").Value = 27.027027027027
$ws.Range("M18").Value = -53
$ws.Range("N18").Value = -88.494492044063
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 35.714285714285
$ws.Range("I19").Value = 465
$ws.Range("J19").Value = 332
$ws.Range("K19").Value = 40.060240963855
$ws.Range("L19").Value = 37.982195845697
$ws.Range("M19").Value = 111.363636363636
$ws.Range("N19").Value = 49.517684887459
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 42
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 10.526315789473
$ws.Range("I20").Value = 248
$ws.Range("J20").Value = 280
$ws.Range("K20").Value = -11.428571428571
$ws.Range("L20").Value = 32.620320855615
$ws.Range("M20").Value = 92.248062015503
$ws.Range("N20").Value = -74.898785425101
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 34.375
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = 10.714285714285
$ws.Range("I21").Value = 1219
$ws.Range("J21").Value = 1118
$ws.Range("K21").Value = 9.033989266547
$ws.Range("L21").Value = 34.547461368653
$ws.Range("M21").Value = 43.243243243243
$ws.Range("N21").Value = -54.531891085415
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 11.111111111111
$ws.Range("L22").Value = 11.111111111111
$ws.Range("M22").Value = -16.666666666666
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 75
$ws.Range("K23").Value = -17.333333333333
$ws.Range("L23").Value = 10.714285714285
$ws.Range("M23").Value = 63.157894736842
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -31.578947368421
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -27.826086956521
$ws.Range("I24").Value = 733
$ws.Range("J24").Value = 833
$ws.Range("K24").Value = -12.004801920768
$ws.Range("L24").Value = 5.164992826398
$ws.Range("M24").Value = 56.289978678038
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -61.111111111111
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -52.830188679245
$ws.Range("I25").Value = 290
$ws.Range("J25").Value = 347
$ws.Range("K25").Value = -16.426512968299
$ws.Range("L25").Value = 3.202846975088
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = 31.111111111111
$ws.Range("I26").Value = 316
$ws.Range("J26").Value = 304
$ws.Range("K26").Value = 3.947368421052
$ws.Range("L26").Value = 6.397306397306
$ws.Range("M26").Value = -6.784660766961
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = 10
$ws.Range("L27").Value = 4.761904761904
$ws.Range("D28").Value = 3
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -70
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = -28.888888888888
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("M29").Value = -35.294117647058
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("M30").Value = -23.076923076923
